# Update the worksheet date and all the three-digit x one-digit
# multiplication problems/answers to the new values.

$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-10 Monday", "2024-06-11 Tuesday"),
    @("738×5=3690", "197×9=1773"),
    @("518×6=3108", "948×7=6636"),
    @("461×3=1383", "227×2=454"),
    @("501×2=1002", "371×7=2597"),
    @("179×4=716", "917×4=3668"),
    @("405×4=1620", "915×6=5490"),
    @("615×3=1845", "230×8=1840"),
    @("266×7=1862", "682×8=5456"),
    @("922×7=6454", "383×5=1915"),
    @("522×6=3132", "232×4=928"),
    @("527×6=3162", "264×8=2112"),
    @("171×5=855", "451×8=3608"),
    @("269×3=807", "342×6=2052"),
    @("140×4=560", "857×2=1714"),
    @("340×2=680", "484×6=2904"),
    @("965×8=7720", "405×2=810"),
    @("233×2=466", "255×3=765"),
    @("950×9=8550", "451×3=1353"),
    @("904×3=2712", "272×6=1632"),
    @("445×5=2225", "223×9=2007"),
    @("285×2=570", "313×9=2817"),
    @("976×9=8784", "974×2=1948"),
    @("184×8=1472", "624×7=4368"),
    @("926×9=8334", "540×9=4860"),
    @("652×9=5868", "933×8=7464")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found for replacement: '$old' -> '$new'"
    }
}

Write-Host "Done applying replacements."
